$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.313.51'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.791.55'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.51'
$ws.Range("E5").Value = '  -0.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5366'
$ws.Range("E7").Value = '  -1.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3763'
$ws.Range("E8").Value = '  -1.74%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07482'
$ws.Range("E9").Value = '  -1.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.31'
$ws.Range("E10").Value = '  -3.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.094'
$ws.Range("E11").Value = '  -2.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.45'
$ws.Range("E13").Value = '  -3.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.096'
$ws.Range("E14").Value = '  -1.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.256'
$ws.Range("E15").Value = '  -0.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.786.00'
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.16'
$ws.Range("E17").Value = '  -2.53%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001057'
$ws.Range("E18").Value = '  -1.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06493'
$ws.Range("E19").Value = '  +0.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.37'
$ws.Range("E21").Value = '  +0.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.935'
$ws.Range("E22").Value = '  -0.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.332.83'
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.10'
$ws.Range("E24").Value = '  -1.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.091'
$ws.Range("E25").Value = '  -4.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.45'
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.27'
$ws.Range("E27").Value = '  -1.65%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.990.94'
$ws.Range("E28").Value = '  -1.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.297'
$ws.Range("E29").Value = '  -6.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.80'
$ws.Range("E30").Value = '  -1.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.093'
$ws.Range("E31").Value = '  -4.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1047'
$ws.Range("E32").Value = '  +3.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.663'
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.531'
$ws.Range("E34").Value = '  -3.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2256'
$ws.Range("E35").Value = '  -2.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06445'
$ws.Range("E36").Value = '  +3.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02276'
$ws.Range("E37").Value = '  -1.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.001'
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.507'
$ws.Range("E39").Value = '  -4.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6156'
$ws.Range("E40").Value = '  -3.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.184'
$ws.Range("E41").Value = '  +2.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.437'
$ws.Range("E42").Value = '  +3.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.03'
$ws.Range("E43").Value = '  -4.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.38'
$ws.Range("E45").Value = '  -0.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.665'
$ws.Range("E46").Value = '  -0.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5763'
$ws.Range("E47").Value = '  -3.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.24'
$ws.Range("E48").Value = '  +0.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.198'
$ws.Range("E49").Value = '  +4.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.930'
$ws.Range("E50").Value = '  -2.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06849'
$ws.Range("E51").Value = '  -0.83%  '
